$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 147, pushing existing rows 147-161 down to 148-162
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with the new weekly record
$ws.Cells.Item(147, 1).Value = 7
$ws.Cells.Item(147, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(147, 3).Value = "Ñuble"
$ws.Cells.Item(147, 4).Value = 45124
$ws.Cells.Item(147, 5).Value = 16
$ws.Cells.Item(147, 6).Value = "Fruta"
$ws.Cells.Item(147, 7).Value = 100108
$ws.Cells.Item(147, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(147, 9).Value = 100108002
$ws.Cells.Item(147, 10).Value = "Mango"
$ws.Cells.Item(147, 11).Value = "Sin especificar"
$ws.Cells.Item(147, 12).Value = "Primera"
$ws.Cells.Item(147, 13).Value = 40
$ws.Cells.Item(147, 14).Value = 9000
$ws.Cells.Item(147, 15).Value = 9000
$ws.Cells.Item(147, 16).Value = 9000
$ws.Cells.Item(147, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(147, 18).Value = "Brasil"
$ws.Cells.Item(147, 19).Value = 2250
$ws.Cells.Item(147, 20).Value = 4
